$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of data (row 8), following the same pattern as existing rows
$row = 8

$ws.Cells.Item($row, 1).Value = 8
$ws.Cells.Item($row, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item($row, 3).Value = "Coquimbo"
$ws.Cells.Item($row, 4).Value = 44595
$ws.Cells.Item($row, 4).NumberFormat = $ws.Cells.Item(2, 4).NumberFormat
$ws.Cells.Item($row, 5).Value = 4
$ws.Cells.Item($row, 6).Value = "Fruta"
$ws.Cells.Item($row, 7).Value = 100107
$ws.Cells.Item($row, 8).Value = "Otros"
$ws.Cells.Item($row, 9).Value = 100107011
$ws.Cells.Item($row, 10).Value = "Tuna"
$ws.Cells.Item($row, 11).Value = "Sin especificar"
$ws.Cells.Item($row, 12).Value = "Primera"
$ws.Cells.Item($row, 13).Value = 200
$ws.Cells.Item($row, 14).Value = 15500
$ws.Cells.Item($row, 15).Value = 16000
$ws.Cells.Item($row, 16).Value = 15750
$ws.Cells.Item($row, 17).Value = "`$/caja 18 kilos"
$ws.Cells.Item($row, 18).Value = "Provincia de Limarí"
$ws.Cells.Item($row, 19).Value = 875
$ws.Cells.Item($row, 20).Value = 18
